# Add season-record columns (Wins / Losses / Ties) to the KCR_2023 sheet.
#
# The previous scrape only pulled team statistics, not the season record.
# This adds three new trailing columns (AD:AF) with the team's win/loss/tie
# totals for the season, repeated on every player row, plus a header row
# matching the existing header formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, formatted exactly like the existing header row (bold,
# centered, top-aligned, thin border) by copying the format from the last
# existing header cell (AC1) before writing the new labels.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record for every player row: 56 wins, 106 losses, 0 ties.
$lastRow = 60
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 56   # AD
    $ws.Cells.Item($r, 31).Value = 106  # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
